# Edit LOQ4259.xlsx per commit diff:
#  - Remove the standalone "5840560 - Marco Antonio Carvalho Pereira" row
#    (old row 13), shifting subsequent rows up by one.
#  - Update the B/C values of several rows with their new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 13 (the "5840560 - Marco Antonio Carvalho Pereira" row with
# no label in column A). This shifts rows 14-24 up to become rows 13-23 and
# updates the sheet dimension to A1:C23 automatically.
$ws.Rows.Item(13).Delete()

# Now fix up the cells whose text content changed relative to a pure shift.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."

$ws.Range("B20").Value = "Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)"
$ws.Range("C20").Value = "Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)"

$ws.Range("B21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação."
